$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 5 (which currently holds "DOC") to make
# room for the new "TC" and "TN" parameter rows.
$ws.Rows("5:6").Insert()

# Populate the new row 5 ("TC") and row 6 ("TN").
# Write the combined ("C") column text first, then the simple name
# columns, so the new shared-string entries land in the same order as
# the target workbook (TC (mg/L), TN (mg/L), TC, TN).
$ws.Range("C5").Value = "TC (mg/L)"
$ws.Range("C6").Value = "TN (mg/L)"

$ws.Range("A5").Value = "TC"
$ws.Range("D5").Value = "TC"

$ws.Range("A6").Value = "TN"
$ws.Range("D6").Value = "TN"

$ws.Range("B5").Value = "(mg/L)"
$ws.Range("B6").Value = "(mg/L)"

# Nudge the small floating textbox shape back down to its anchor cell
# now that two rows were inserted above it (row insert does not retarget
# floating shape anchors automatically in this engine). The two new rows
# (5 and 6) use the default 16pt row height, so the shape moves down by
# 32pt to land back on the same relative offset within its anchor row.
$shp = $ws.Shapes.Item(1)
$shp.Top = 232.75

# Update the active selection to match the committed workbook.
[void]$ws.Range("A7").Select()
